# Generate Report for Handoff
# Re-order the localization-status rows so the b1a02e00 entry (now ready for
# a new handoff) moves to the bottom, the ffff95b24df0 entry moves up to
# row 2, and the ffffff16b0a95d entry moves up to row 3. Update the status,
# dates and error detail for the b1a02e00 entry to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff95b24df0-adc0-4391-8852-e6f51cbcdd72.md"
$ws1.Range("B2").Value = "e2e\ffff95b24df0-adc0-4391-8852-e6f51cbcdd72.md"
$ws1.Range("G2").Value = "2016-09-01 05:07:14"

$ws1.Range("A3").Value = "ffffff16b0a95d-ee72-4fcc-b4d9-2b8fd9889e2c.md"
$ws1.Range("B3").Value = "e2e\ffffff16b0a95d-ee72-4fcc-b4d9-2b8fd9889e2c.md"

$ws1.Range("A4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws1.Range("B4").Value = "e2e\b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-09-01 05:10:21"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff95b24df0-adc0-4391-8852-e6f51cbcdd72.md"
$ws2.Range("G2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-01 05:07:07"
$ws2.Range("I2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.md"
$ws2.Range("J2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-01 05:07:33"

$ws2.Range("A3").Value = "ffffff16b0a95d-ee72-4fcc-b4d9-2b8fd9889e2c.md"
$ws2.Range("F3").Value = "True"

$ws2.Range("A4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("F4").Value = "False"
$ws2.Range("G4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.fafa3fe2a12a5a6c5a4f300bf93b291105c68f7f.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-09-01 05:10:17"
$ws2.Range("I4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws2.Range("J4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.fafa3fe2a12a5a6c5a4f300bf93b291105c68f7f.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-09-01 05:09:41"
$ws2.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f94f3ddcde69aefca00a304d3106e6d189217b1/e2e/b1a02e00-2cc8-4dc7-9978-312c489ae804.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0c68e82f89ab85159042a7b862264de812725e/e2e/b1a02e00-2cc8-4dc7-9978-312c489ae804.md."

$ws2.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff95b24df0-adc0-4391-8852-e6f51cbcdd72.md"
$ws3.Range("G2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-01 05:07:14"
$ws3.Range("I2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.md"
$ws3.Range("J2").Value = "e03e8cad-5f2b-4759-a80d-0581a2aa856c.101179c73998b821a8504f720cbefac42762ec1d.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-01 05:07:41"

$ws3.Range("A3").Value = "ffffff16b0a95d-ee72-4fcc-b4d9-2b8fd9889e2c.md"
$ws3.Range("F3").Value = "True"

$ws3.Range("A4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("F4").Value = "False"
$ws3.Range("G4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.fafa3fe2a12a5a6c5a4f300bf93b291105c68f7f.de-de.xlf"
$ws3.Range("H4").Value = "2016-09-01 05:10:21"
$ws3.Range("I4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.md"
$ws3.Range("J4").Value = "b1a02e00-2cc8-4dc7-9978-312c489ae804.fafa3fe2a12a5a6c5a4f300bf93b291105c68f7f.de-de.xlf"
$ws3.Range("K4").Value = "2016-09-01 05:09:48"
$ws3.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f94f3ddcde69aefca00a304d3106e6d189217b1/e2e/b1a02e00-2cc8-4dc7-9978-312c489ae804.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0c68e82f89ab85159042a7b862264de812725e/e2e/b1a02e00-2cc8-4dc7-9978-312c489ae804.md."

$ws3.Columns.Item(16).ColumnWidth = 39.15

$wb.Save()
